$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha Ponto")

# Row 6 ("Arrumando teleporte dos tubos brancos e pretos") ran 10 minutes
# longer than before: end time moves from 11:30 to 11:45 (fixing the fade
# bug took a bit longer).
$ws.Range("D6").Value = 0.48958333333333331

# Copy the date/time number formats from row 6 down into the new row 7
# before filling in values, so the new cells keep the same styling
# (date format on B, time format on C/D) instead of Excel's generic
# auto-format.
$ws.Range("B6:D6").Copy()
$ws.Range("B7:D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)

# New entry: 11:45 - 12:15 working on SITS, "Fixed fade bug".
$ws.Range("B7").Value = 44902
$ws.Range("C7").Value = 0.48958333333333331
$ws.Range("D7").Value = 0.51041666666666663
$ws.Range("F7").Value = "SITS"
$ws.Range("G7").Value = "Fixed fade bug"

# Active cell moves on to the next row's end-time entry.
$ws.Range("D8").Select()

$wb.Save()
